# C1--C2-and-C3-PowerPoint.pptx edit
# 1) Slide 16 (the cash-flow "PLENARY" table) gets a new table style.
# 2) The deck's theme palette switches from the "Integral" colours to the
#    stock "Office Theme" colours (Design swap), applied through the
#    presentation's theme colour scheme.

$p = $ppt.ActivePresentation

# --- 1. Swap the table style on the slide-16 table -------------------------
$tableSlide  = $p.Slides.Item(16)
$tableShape  = $tableSlide.Shapes.Item(3)
$tableShape.Table.ApplyStyle("{DA1C4308-F5E8-4203-AA03-5D38386E6B8E}")

# --- 2. Re-colour the theme: Integral -> Office Theme -----------------------
$colors = $p.SlideMaster.Theme.ThemeColorScheme

$colors.Item(1).RGB  = 0x000000   # dk1
$colors.Item(2).RGB  = 0xFFFFFF   # lt1
$colors.Item(3).RGB  = 0x6A5444   # dk2      (44546A)
$colors.Item(4).RGB  = 0xE6E6E7   # lt2      (E7E6E6)
$colors.Item(5).RGB  = 0xD59B5B   # accent1  (5B9BD5)
$colors.Item(6).RGB  = 0x317DED   # accent2  (ED7D31)
$colors.Item(7).RGB  = 0xA5A5A5   # accent3  (A5A5A5)
$colors.Item(8).RGB  = 0x00C0FF   # accent4  (FFC000)
$colors.Item(9).RGB  = 0xC47244   # accent5  (4472C4)
$colors.Item(10).RGB = 0x47AD70   # accent6  (70AD47)
$colors.Item(11).RGB = 0xC16305   # hlink    (0563C1)
$colors.Item(12).RGB = 0x724F95   # folHlink (954F72)
